# "final step of my project"
# The lat/long columns (B:C) were converted from text-formatted shared
# strings ("50.450100", "30.523400", ...) into real numeric values with
# full floating point precision, and the active selection moved to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 50.447730999999997
$ws.Range("C2").Value = 30.542721
$ws.Range("B3").Value = 49.839683999999998
$ws.Range("C3").Value = 24.029716000000001
$ws.Range("B4").Value = 48.922634000000002
$ws.Range("C4").Value = 24.711117000000002

# Column widths were nudged slightly (new Excel install / font metrics).
$ws.Columns.Item(1).ColumnWidth = 17.6
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(3).ColumnWidth = 11.3

# Selection moved from C11 to C3.
[void]$ws.Range("C3").Select()
